# Remove the two rows describing the "Immunocompetent Mouse Model for
# Crimean-Congo Hemorrhagic Fever Virus" entry (RefID 13) and its
# companion "Direct Submission" entry (RefID 14) that shared the same
# accession numbers (MW058028, MW058029, MW058030).
#
# These correspond to worksheet rows 6 and 7 (the header is row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Select rows 6:7 (mirrors what a user would do before deleting them)
# and delete the entire rows, shifting everything below up.
$rng = $ws.Range("A6:A7").EntireRow
$rng.Select()
$rng.Delete()

$wb.Save()
